$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.236.18'
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('D3').Value = '1.829.94'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.10'
$ws.Range('E5').Value = '  -1.68%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6019'
$ws.Range('E6').Value = '  -2.93%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.004'
$ws.Range('E7').Value = '  +0.26%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.06951'
$ws.Range('E8').Value = '  -5.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2766'
$ws.Range('E9').Value = '  -3.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.56'
$ws.Range('E10').Value = '  -4.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07609'
$ws.Range('E11').Value = '  -1.35%  '
$ws.Range('D12').Value = '1.834.70'
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.750'
$ws.Range('E13').Value = '  -3.74%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6285'
$ws.Range('E14').Value = '  -4.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000009850'
$ws.Range('E15').Value = '  -6.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '77.71'
$ws.Range('E16').Value = '  -4.40%  '
$ws.Range('D17').Value = '28.818.89'
$ws.Range('E17').Value = '  -1.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.573'
$ws.Range('E18').Value = '  -10.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '216.99'
$ws.Range('E19').Value = '  -8.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.004'
$ws.Range('E20').Value = '  +0.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.56'
$ws.Range('E21').Value = '  -4.96%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.889'
$ws.Range('E22').Value = '  -4.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.002'
$ws.Range('E23').Value = '  -0.50%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '155.81'
$ws.Range('E24').Value = '  -0.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.960'
$ws.Range('E25').Value = '  -5.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1290'
$ws.Range('E26').Value = '  -2.90%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.50'
$ws.Range('E27').Value = '  -4.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06441'
$ws.Range('E28').Value = '  -6.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.420'
$ws.Range('E29').Value = '  -3.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.442'
$ws.Range('E30').Value = '  -2.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.836'
$ws.Range('E31').Value = '  -1.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.784'
$ws.Range('E32').Value = '  -5.41%  '
$ws.Range('B33').Value = 'LidoDAOToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.728'
$ws.Range('E33').Value = '  -0.58%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.095'
$ws.Range('E34').Value = '  -4.74%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6477'
$ws.Range('E35').Value = '  -4.68%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.539'
$ws.Range('E36').Value = '  -1.66%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.752'
$ws.Range('E37').Value = '  -0.93%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01758'
$ws.Range('E38').Value = '  -3.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.512'
$ws.Range('E39').Value = '  -2.11%  '
$ws.Range('D40').Value = '1.143.69'
$ws.Range('E40').Value = '  -6.96%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8906'
$ws.Range('E41').Value = '  -5.60%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.003'
$ws.Range('E42').Value = '  +0.21%  '
$ws.Range('D43').Value = '1.989.23'
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.83'
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '62.07'
$ws.Range('E45').Value = '  -4.29%  '
$ws.Range('E46').Value = '  -4.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.617'
$ws.Range('E47').Value = '  -3.75%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.499'
$ws.Range('E48').Value = '  -2.99%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4548'
$ws.Range('E49').Value = '  -0.45%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05494'
$ws.Range('E50').Value = '  -2.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.412'
$ws.Range('E51').Value = '  -6.59%  '
